# Updates the Price (column D) and Volume(1h) (column E) columns of the
# cryptos list with refreshed values, matching a new data snapshot.
# Values in column D that look like plain numbers are written with a
# leading apostrophe ('') so Excel keeps storing them as text (matching
# the original inlineStr/text storage) instead of silently converting
# them to floating point numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.410.07'
$ws.Range('E2').Value = '  +0.44%  '
$ws.Range('D3').Value = '1.870.20'
$ws.Range('E3').Value = '  +0.05%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').Value = '''330.42'
$ws.Range('E5').Value = '  -2.84%  '
$ws.Range('D6').Value = '''1.001'
$ws.Range('E6').Value = '  +0.01%  '
$ws.Range('D7').Value = '''0.4616'
$ws.Range('E7').Value = '  -1.89%  '
$ws.Range('D8').Value = '''0.4010'
$ws.Range('E8').Value = '  +1.80%  '
$ws.Range('D9').Value = '''47.81'
$ws.Range('E9').Value = '  +1.16%  '
$ws.Range('D10').Value = '''0.07858'
$ws.Range('E10').Value = '  -1.85%  '
$ws.Range('D11').Value = '''0.9858'
$ws.Range('E11').Value = '  -1.97%  '
$ws.Range('D12').Value = '''21.32'
$ws.Range('E12').Value = '  -2.62%  '
$ws.Range('D13').Value = '1.869.14'
$ws.Range('E13').Value = '  -0.55%  '
$ws.Range('D14').Value = '''5.850'
$ws.Range('E14').Value = '  -2.57%  '
$ws.Range('D15').Value = '''6.993'
$ws.Range('E15').Value = '  -4.01%  '
$ws.Range('D16').Value = '''1.001'
$ws.Range('E16').Value = '  -0.12%  '
$ws.Range('D17').Value = '''88.18'
$ws.Range('E17').Value = '  -3.34%  '
$ws.Range('D18').Value = '''0.06534'
$ws.Range('E18').Value = '  -1.25%  '
$ws.Range('D19').Value = '''0.00001019'
$ws.Range('E19').Value = '  -2.42%  '
$ws.Range('D20').Value = '''17.23'
$ws.Range('E20').Value = '  -2.37%  '
$ws.Range('D21').Value = '''0.9988'
$ws.Range('E21').Value = '  -0.25%  '
$ws.Range('D22').Value = '28.388.34'
$ws.Range('E22').Value = '  +0.34%  '
$ws.Range('D23').Value = '''5.346'
$ws.Range('E23').Value = '  -1.95%  '
$ws.Range('E24').Value = '  -1.83%  '
$ws.Range('E25').Value = '  -1.78%  '
$ws.Range('D26').Value = '2.090.88'
$ws.Range('E26').Value = '  -0.50%  '
$ws.Range('D27').Value = '''157.44'
$ws.Range('E27').Value = '  -1.57%  '
$ws.Range('E28').Value = '  -2.28%  '
$ws.Range('D29').Value = '''2.059'
$ws.Range('E29').Value = '  -4.07%  '
$ws.Range('D30').Value = '''5.294'
$ws.Range('E30').Value = '  -3.94%  '
$ws.Range('D31').Value = '''117.52'
$ws.Range('E31').Value = '  -2.34%  '
$ws.Range('D32').Value = '''0.9578'
$ws.Range('E32').Value = '  -2.10%  '
$ws.Range('D33').Value = '''0.09341'
$ws.Range('E33').Value = '  -1.85%  '
$ws.Range('D34').Value = '''3.583'
$ws.Range('E34').Value = '  -0.29%  '
$ws.Range('D35').Value = '''1.387'
$ws.Range('E35').Value = '  +0.58%  '
$ws.Range('D36').Value = '''5.245'
$ws.Range('E36').Value = '  -2.09%  '
$ws.Range('D37').Value = '''0.06031'
$ws.Range('E37').Value = '  -1.09%  '
$ws.Range('D38').Value = '''0.02201'
$ws.Range('E38').Value = '  -3.34%  '
$ws.Range('D39').Value = '''8.293'
$ws.Range('E40').Value = '  -1.82%  '
$ws.Range('D41').Value = '''0.9997'
$ws.Range('D42').Value = '''0.5757'
$ws.Range('E42').Value = '  -3.74%  '
$ws.Range('D43').Value = '''0.1809'
$ws.Range('E43').Value = '  -3.91%  '
$ws.Range('D44').Value = '''10.04'
$ws.Range('E44').Value = '  -3.43%  '
$ws.Range('D45').Value = '''1.266'
$ws.Range('E45').Value = '  -1.78%  '
$ws.Range('D46').Value = '''2.289'
$ws.Range('E46').Value = '  +12.44%  '
$ws.Range('D47').Value = '''0.5435'
$ws.Range('E47').Value = '  -3.25%  '
$ws.Range('D48').Value = '''11.90'
$ws.Range('E48').Value = '  -1.92%  '
$ws.Range('D49').Value = '''0.07124'
$ws.Range('E49').Value = '  +3.14%  '
$ws.Range('D50').Value = '''1.887'
$ws.Range('D51').Value = '''111.34'
$ws.Range('E51').Value = '  +0.01%  '
